# Add implementation for version 2
#
# The "мини-печь StarWind" slide (slide 2) had its title placeholder
# ("Title 1" / "мини-печь StarWind") removed entirely from the slide.
#
# NOTE: calling .Delete() on a layout placeholder shape only resets it
# (clears its text / reassigns a fresh id) instead of removing the <p:sp>
# element from the slide's shape tree - mirroring how PowerPoint keeps
# placeholder "ghosts" around. Using .Cut() instead removes the shape
# from the slide completely, which is what the target diff requires.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$title = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "мини-печь StarWind") {
            $title = $shape
            break
        }
    }
}

if ($null -eq $title) {
    # Fall back to the title placeholder of the slide, in case the text
    # match above ever fails to resolve.
    $title = $s.Shapes.Placeholders.Item(1)
}

$title.Cut()
